# ajustes da reprovacao do caso
#
# 1) " II – As partes têm interesse na rescisão antecipada do Acordo;"
#      -> "II – As partes têm interesse na rescisão antecipada do Acordo."
# 2) "As partes acima qualificadas, ..." -> "As partes, acima qualificadas, ..."
# 3) Remove ", na presença das testemunhas abaixo assinadas" from the closing
#    paragraph.
# 4) signHere anchors (school e-mail / worker e-mail) get orange highlight +
#    white text and their spacing collapses to 0/0; the underline paragraphs
#    right after them also collapse their spacing to 0/0.

$d = $word.ActiveDocument

# --- 1) Clause II punctuation / leading space -------------------------------
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "II – As partes têm interesse na rescisão antecipada do Acordo."

# --- 2) add comma after "As partes" ----------------------------------------
$p9 = $d.Paragraphs.Item(9)
$p9start = $p9.Range.Start
$oldLead = "As partes acima qualificadas, decidem acordar o presente Termo de Rescisão do "
$newLead = "As partes, acima qualificadas, decidem acordar o presente Termo de Rescisão do "
$leadRange = $d.Range($p9start, $p9start + $oldLead.Length)
$leadRange.Text = $newLead

# --- 3) drop the witnesses clause -------------------------------------------
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Text = "Por estarem assim justas e acertadas, as partes firmam o presente documento em 2 (duas) vias de igual teor e forma, para que produza seus efeitos legais."

# --- 4) signHere anchors: highlight + spacing -------------------------------
function Set-SignHereAnchor($paragraph) {
    $paragraph.Format.SpaceAfter = 0
    $paragraph.Format.SpaceBefore = 0
    $range = $paragraph.Range
    $range.Font.Color = 16777215           # white (wdColorWhite)
    $range.Shading.Texture = 0             # wdTextureNone -> w:val="clear"
    $range.Shading.BackgroundPatternColor = 39423   # ff9900 (orange)
}

$p17 = $d.Paragraphs.Item(17)   # {{ generate_anchor('signHere', school["email"]) }}
Set-SignHereAnchor $p17

$p18 = $d.Paragraphs.Item(18)   # underline row right below it
$p18.Format.SpaceAfter = 0
$p18.Format.SpaceBefore = 0

$p22 = $d.Paragraphs.Item(22)   # {{ generate_anchor('signHere', item.email) }}
Set-SignHereAnchor $p22

$p23 = $d.Paragraphs.Item(23)   # underline row right below it
$p23.Format.SpaceAfter = 0
$p23.Format.SpaceBefore = 0

Write-Output "done"
